$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$names = @("ContextualSpacing","NoSpaceBetweenParagraphsOfSameStyle","DoNotAddSpaceBetweenParagraphsOfSameStyle","SpaceBetweenParagraphsOfSameStyle","AddSpaceBetweenParagraphsOfSameStyle")
foreach ($n in $names) {
    try {
        $v = $p1.Format.$n
        Write-Output "GET Format.$n = $v"
    } catch {
        Write-Output "GET Format.$n error: $_"
    }
    try {
        $v = $p1.$n
        Write-Output "GET Paragraph.$n = $v"
    } catch {
        Write-Output "GET Paragraph.$n error: $_"
    }
}
